$d = $word.ActiveDocument

# Word stores each paragraph's formatting (w:pPr, paraId, ...) on the paragraph
# MARK that sits at the very end of the paragraph's own text. In the target
# revision, a new "Sample plots for the two events" bullet is inserted right
# before three existing bullets; the ORIGINAL paragraph mark (and its paraId)
# ends up on the new "Sample plots..." text, while the pre-existing text is
# pushed into a freshly created paragraph mark with no paraId. To reproduce
# that, split the paragraph right before its own trailing mark (which keeps
# the original mark, and its paraId, attached to the first/old-text part),
# then swap: the first part becomes "Sample plots for the two events" and the
# new, attribute-less second part receives the original text.

function Split-BulletAndSwap($paraIndex, $expectedPrefix, $newFirstText, $origLevel) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $origText = $r.Text.Substring(0, $r.Text.Length - 1)

    if (-not $origText.StartsWith($expectedPrefix)) {
        throw "Paragraph $paraIndex does not start with expected text '$expectedPrefix' (found: '$origText')"
    }

    $ins = $d.Range($r.End - 1, $r.End - 1)
    $ins.InsertParagraphAfter()

    # First paragraph (keeps the original paraId): scope Find to just this
    # paragraph so later, identical text elsewhere isn't touched instead.
    $p1 = $d.Paragraphs($paraIndex)
    $p1.Range.Find.Execute($origText, $true, $false, $false, $false, $false, `
        $true, 1, $false, $newFirstText, 2) | Out-Null
    $p1 = $d.Paragraphs($paraIndex)
    $p1.Range.ListFormat.ListLevelNumber = 2   # ilvl=1 ("Sample plots..." is a sub-bullet)

    # Second paragraph (new, no paraId): gets the original text back.
    $p2 = $d.Paragraphs($paraIndex + 1)
    $p2.Range.ListFormat.ListLevelNumber = $origLevel + 1
    $p2.Range.Text = $origText
}

# Process from the bottom of the document upward so that earlier paragraph
# indices remain valid while later ones are being edited.

# Section "Three IBRs": insert "Sample plots..." before "IBRs are modeled..."
Split-BulletAndSwap 13 "IBR" "Sample plots for the two events" 1

# Section "Single IBR": insert "Sample plots..." before "IBR is modeled..."
Split-BulletAndSwap 10 "IBR is modeled" "Sample plots for the two events" 1

# Start of "Single IBR" section: insert "Sample plots..." before
# "Enhanced IEEE 39-Bus System_Single IBR"
Split-BulletAndSwap 8 "Enhanced IEEE 39-Bus" "Sample plots for the two events" 0

# Collapse the remaining multi-run / spell-checked paragraphs into single
# runs by replacing their (unchanged) concatenated text via Find & Replace --
# this merges the text into a single run and drops the proofErr spell-check
# markers, matching the target revision.

$d.Content.Find.Execute(
    "Enhanced IEEE 39-Bus System_Three IBRs", $true, $false, $false, $false, $false,
    $true, 1, $false, "Enhanced IEEE 39-Bus System_Three IBRs", 2) | Out-Null

$d.Content.Find.Execute(
    "Contains PSSE, PSLF, PSCAD files of improved IEEE 39 Bus System with three identical IBRs for running dynamic simulations with two events- Bus fault at bus 16, Generator trip at Gen 32",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Contains PSSE, PSLF, PSCAD files of improved IEEE 39 Bus System with three identical IBRs for running dynamic simulations with two events- Bus fault at bus 16, Generator trip at Gen 32",
    2) | Out-Null

$d.Content.Find.Execute(
    "Contains PSSE, PSLF, PSCAD files of improved IEEE 39 Bus System with single IBR for running dynamic simulations with two events- Bus fault at bus 16, Generator trip at Gen 32",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Contains PSSE, PSLF, PSCAD files of improved IEEE 39 Bus System with single IBR for running dynamic simulations with two events- Bus fault at bus 16, Generator trip at Gen 32",
    2) | Out-Null
